# RDCC-3540 Upload Template Improvements
$wb = $excel.ActiveWorkbook

# Rename the first sheet from "Case Worker Data" to "Staff Data"
$ws = $wb.Worksheets.Item("Case Worker Data")
$ws.Name = "Staff Data"

# Add two new header columns first (so their shared strings are registered before the Service1-8 ones)
$ws.Range("V1").Value = "Task Supervisor"
$ws.Range("W1").Value = "Case Allocator"

# Update Area of WorkN headers (L1:S1) to ServiceN
$ws.Range("L1").Value = "Service1"
$ws.Range("M1").Value = "Service2"
$ws.Range("N1").Value = "Service3"
$ws.Range("O1").Value = "Service4"
$ws.Range("P1").Value = "Service5"
$ws.Range("Q1").Value = "Service6"
$ws.Range("R1").Value = "Service7"
$ws.Range("S1").Value = "Service8"

# Leave selection on the last edited header cell
$null = $ws.Range("W1").Select()
